# fix(pipelining): fix testing programs #18
#
# The LOAD instructions in this pipelining example are rewritten to use
# post-indexed addressing (no write-back hazard on the same cycle as the
# load), and the MOV that initializes R0 is adjusted accordingly. This
# pushes the "ST" (store/write-back) pipeline stage one column later for
# each subsequent instruction in the "Pipelining" diagram sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Code"
$ws2 = $wb.Worksheets.Item(2)   # "Pipelining"

# --- Sheet "Code": update instruction text -------------------------------
# Update the LOAD instructions first (so new shared strings for them are
# created before the MOV string), matching the order the strings end up
# added to the workbook.
$ws1.Range("B10").Value = "LOAD R2, [R0], #1"
$ws1.Range("B11").Value = "LOAD R3, [R0], #2"
$ws1.Range("B8").Value  = "MOV R0, #0 ROR 0"

# --- Sheet "Pipelining": update instruction text and shift ST/F/E1/E2 ----
$ws2.Range("B3").Value = "LOAD R2, [R0], #1"
$ws2.Range("B4").Value = "LOAD R3, [R0], #2"

# Row 4 pipeline stages shift one column to the right (G:J -> H:K)
$ws2.Range("G4:J4").ClearContents()
$ws2.Range("H4").Value = "ST"
$ws2.Range("I4").Value = "F"
$ws2.Range("J4").Value = "E1"
$ws2.Range("K4").Value = "E2"

# Row 5 pipeline stages shift two columns to the right (H:K -> J:M)
$ws2.Range("H5:K5").ClearContents()
$ws2.Range("J5").Value = "ST"
$ws2.Range("K5").Value = "F"
$ws2.Range("L5").Value = "E1"
$ws2.Range("M5").Value = "E2"

# --- View state (best effort) --------------------------------------------
$ws1.Range("B10").Select()
$ws2.Range("E2").Select()
$ws2.Activate()
